$d = $word.ActiveDocument

# 1) Title: "MTN mPulse GAP" -- collapse the three runs ("MTN ", "mPulse", " GAP")
#    into a single run by replacing the phrase with itself.
$d.Content.Find.Execute(
    "MTN mPulse GAP",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "MTN mPulse GAP",
    2
)

# 2) Kagukar paragraph: collapse "The " + "Kagukar" + " Mentorship track kicked off
#    with the introductory session" into a single run, then extend that paragraph
#    with the additional sentence describing the session.
$d.Content.Find.Execute(
    "The Kagukar Mentorship track kicked off with the introductory session",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The Kagukar Mentorship track kicked off with the introductory session",
    2
)

$r = $d.Content
$r.Find.Execute(
    "The Kagukar Mentorship track kicked off with the introductory session",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
)
$r.Collapse(0)
$r.InsertAfter(" where a TechQuest Representative conducted a formal introduction between the various mentorship tracks, the Mentees and the corresponding Mentors respectively, the session lasted for 45mins followed by a breakout session --with one Mentee on the Kagukar Web Application track-- wherein the Mentee had a One-One introduction with the Mentor discussing the whole concept of the application and it's importance.")

# 3) "Joshua Agboola" (appears twice) -- collapse "Joshua " + "Agboola" into a
#    single run in both locations.
$d.Content.Find.Execute(
    "Joshua Agboola",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Joshua Agboola",
    2
)
